$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 158; this shifts the former rows 158..258 down to 159..259,
# preserving all their values/formats.
$ws.Rows(158).Insert()

# The new row 158 should start as a duplicate of the data that used to live there
# (now shifted to row 159), except for a new "Fecha" (date) value.
for ($col = 1; $col -le 18; $col++) {
    $src = $ws.Cells.Item(159, $col)
    $dst = $ws.Cells.Item(158, $col)
    $dst.Value = $src.Value2
}

# New weekly date for the inserted record: 2023-09-07 (serial 45176)
$ws.Cells.Item(158, 4).Value = 45176
